$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-report row was inserted as the new row 98 (Haba /
# Femacal de La Calera), pushing the previous rows 98:162 down to 99:163.
# Use Rows.Insert with xlShiftDown (-4121) so every downstream row keeps
# its original data, exactly like Excel's own "Insert Row" command.
$ws.Rows.Item(98).Insert(-4121)

# Fill in the newly inserted row 98 with its reported values. Columns
# A,B,C,E,F,G,H,I,N,O,Q,R carry the same category/quality/unit/origin
# metadata as the surrounding "Haba" records; D,J,K,L,M,P are the new
# date/volume/price figures for this report.
$ws.Cells.Item(98, 1).Value  = 3
$ws.Cells.Item(98, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(98, 3).Value  = "Coquimbo"
$ws.Cells.Item(98, 4).Value  = 44762
$ws.Cells.Item(98, 5).Value  = 5
$ws.Cells.Item(98, 6).Value  = 100112026
$ws.Cells.Item(98, 7).Value  = "Haba"
$ws.Cells.Item(98, 8).Value  = "Sin especificar"
$ws.Cells.Item(98, 9).Value  = "Primera"
$ws.Cells.Item(98, 10).Value = 73
$ws.Cells.Item(98, 11).Value = 17000
$ws.Cells.Item(98, 12).Value = 18000
$ws.Cells.Item(98, 13).Value = 17479
$ws.Cells.Item(98, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(98, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(98, 16).Value = 699
$ws.Cells.Item(98, 17).Value = 25
$ws.Cells.Item(98, 18).Value = "Hortaliza"
